# "update embassy backend process"
# The embassies header row is renamed:
#   A1: "Embassy of"  -> "Embassy In"
#   B1: "Embassy In"  -> "Embassy City"
# (the old "Embassy of" shared string is dropped and a new "Embassy City"
#  string is introduced; all other header cells C1:H1 are unchanged)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Embassy In"
$ws.Range("B1").Value = "Embassy City"
